# "correct indentation of references"
#
# 1) The document's sole paragraph switches from body-text ("fliesstext")
#    to Heading 1 ("berschrift1").
# 2) The "Heading 1" style ("berschrift1") loses its hard-coded
#    space-before (it keeps only the space-after).
# 3) Two new custom paragraph styles are introduced for bibliography /
#    reference lists, each using a hanging indent so that wrapped lines
#    align under the first word instead of under the number/marker:
#      - "references"    (w:name "references")    -> 19.85pt hanging indent,
#                                                     left-aligned
#      - "references_ni" (w:name "references_ni") -> 14.2pt hanging indent,
#                                                     based on "references"

$d = $word.ActiveDocument

# --- 1) retarget the (only) paragraph onto the Heading 1 style ---------
$d.Paragraphs(1).Style = "berschrift1"

# --- 2) trim the hard-coded space-before on Heading 1 -------------------
$headingStyle = $d.Styles("berschrift1")
$headingStyle.ParagraphFormat.SpaceBefore = 0

# --- 3) add the "references" paragraph style -----------------------------
$references = $d.Styles.Add("references", 1)
$references.BaseStyle = "fliesstext"
$references.NextParagraphStyle = "fliesstext"
$references.QuickStyle = $true
$references.ParagraphFormat.LeftIndent = 19.85
$references.ParagraphFormat.FirstLineIndent = -19.85
$references.ParagraphFormat.Alignment = 0

# --- 3b) add the "references_ni" paragraph style (tighter indent) -------
$referencesNi = $d.Styles.Add("referencesni", 1)
$referencesNi.BaseStyle = "references"
$referencesNi.NextParagraphStyle = "fliesstext"
$referencesNi.QuickStyle = $true
$referencesNi.ParagraphFormat.LeftIndent = 14.2
$referencesNi.ParagraphFormat.FirstLineIndent = -14.2
$referencesNi.NameLocal = "references_ni"
